# Auto-generated Excel COM-interop script
# Applies scheduled market-data / profit recalculation updates to the Spriggan_Profits workbook
# across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1363.1428
$ws.Range("I28").Value = 1363.1428
$ws.Range("K28").Value = 1363.1428
$ws.Range("M28").Value = -878.1428000000001
$ws.Range("H62").Value = 5023.6665
$ws.Range("I62").Value = 4758.143
$ws.Range("K62").Value = 4758.143
$ws.Range("M62").Value = -4134.143
$ws.Range("H65").Value = 5023.6665
$ws.Range("I65").Value = 4758.143
$ws.Range("K65").Value = 23790.715
$ws.Range("M65").Value = -20670.715
$ws.Range("H86").Value = 2695.125
$ws.Range("J86").Value = 2080.3333
$ws.Range("L86").Value = 2080.3333
$ws.Range("N86").Value = -4326.3333
$ws.Range("H89").Value = 2695.125
$ws.Range("J89").Value = 2080.3333
$ws.Range("L89").Value = 10401.6665
$ws.Range("N89").Value = -21633.6665
$ws.Range("H111").Value = 1045.2858
$ws.Range("I111").Value = 852.8333
$ws.Range("J111").Value = 2200
$ws.Range("K111").Value = 2558.4999
$ws.Range("L111").Value = 6600
$ws.Range("M111").Value = 508.5001000000002
$ws.Range("N111").Value = -12734
$ws.Range("H113").Value = 2327.9333
$ws.Range("I113").Value = 2231.8462
$ws.Range("K113").Value = 2231.8462
$ws.Range("M113").Value = 1022.1538
$ws.Range("H137").Value = 1678.4546
$ws.Range("I137").Value = 1333.2632
$ws.Range("K137").Value = 3999.7896
$ws.Range("M137").Value = -1449.7896
$ws.Range("H138").Value = 2610.6135
$ws.Range("I138").Value = 2096.3333
$ws.Range("K138").Value = 6288.999899999999
$ws.Range("M138").Value = -1148.999899999999
$ws.Range("H141").Value = 1950.909
$ws.Range("I141").Value = 1548.579
$ws.Range("J141").Value = 4499
$ws.Range("K141").Value = 4645.737
$ws.Range("L141").Value = 13497
$ws.Range("M141").Value = 534.2629999999999
$ws.Range("N141").Value = -23857
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 644562.9399999999
$ws.Range("I2").Value = 1854093.1
$ws.Range("J2").Value = 4223.4707
$ws.Range("K2").Value = 1854093.1
$ws.Range("L2").Value = 4223.4707
$ws.Range("M2").Value = -1853980.1
$ws.Range("N2").Value = -4449.4707
$ws.Range("H27").Value = 48930.668
$ws.Range("J27").Value = 48930.668
$ws.Range("L27").Value = 48930.668
$ws.Range("N27").Value = -49298.668
$ws.Range("H32").Value = 4679.0957
$ws.Range("I32").Value = 2314.3035
$ws.Range("J32").Value = 12469
$ws.Range("K32").Value = 2314.3035
$ws.Range("L32").Value = 12469
$ws.Range("M32").Value = -2027.3035
$ws.Range("N32").Value = -13043
$ws.Range("H61").Value = 166669380
$ws.Range("I61").Value = 200002240
$ws.Range("K61").Value = 200002240
$ws.Range("M61").Value = -200002028
$ws.Range("H116").Value = 644562.9399999999
$ws.Range("I116").Value = 1854093.1
$ws.Range("J116").Value = 4223.4707
$ws.Range("K116").Value = 1854093.1
$ws.Range("L116").Value = 4223.4707
$ws.Range("M116").Value = -1851799.1
$ws.Range("N116").Value = -8811.4707
$ws.Range("H122").Value = 2065.0527
$ws.Range("I122").Value = 1966.0667
$ws.Range("J122").Value = 2436.25
$ws.Range("K122").Value = 5898.2001
$ws.Range("L122").Value = 7308.75
$ws.Range("M122").Value = -3448.2001
$ws.Range("N122").Value = -12208.75
$ws.Range("H125").Value = 72685.60000000001
$ws.Range("J125").Value = 72685.60000000001
$ws.Range("L125").Value = 72685.60000000001
$ws.Range("N125").Value = -82525.60000000001
$ws.Range("H132").Value = 3036952.8
$ws.Range("I132").Value = 3709842.2
$ws.Range("K132").Value = 11129526.6
$ws.Range("M132").Value = -11126996.6
$ws.Range("H136").Value = 166669380
$ws.Range("I136").Value = 200002240
$ws.Range("K136").Value = 600006720
$ws.Range("M136").Value = -600004170
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 644562.9399999999
$ws.Range("I3").Value = 1854093.1
$ws.Range("J3").Value = 4223.4707
$ws.Range("K3").Value = 1854093.1
$ws.Range("L3").Value = 4223.4707
$ws.Range("M3").Value = -1853979.1
$ws.Range("N3").Value = -4451.4707
$ws.Range("H22").Value = 7938433
$ws.Range("I22").Value = 350
$ws.Range("J22").Value = 15876516
$ws.Range("K22").Value = 350
$ws.Range("L22").Value = 15876516
$ws.Range("M22").Value = -177
$ws.Range("N22").Value = -15876862
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2170.6365
$ws.Range("I22").Value = 2481.5745
$ws.Range("J22").Value = 343.875
$ws.Range("K22").Value = 2481.5745
$ws.Range("L22").Value = 343.875
$ws.Range("M22").Value = -2131.5745
$ws.Range("N22").Value = -1043.875
$ws.Range("H31").Value = 7665.7915
$ws.Range("I31").Value = 2353.0386
$ws.Range("J31").Value = 13944.5
$ws.Range("K31").Value = 2353.0386
$ws.Range("L31").Value = 13944.5
$ws.Range("M31").Value = -2058.0386
$ws.Range("N31").Value = -14534.5
$ws.Range("H34").Value = 7665.7915
$ws.Range("I34").Value = 2353.0386
$ws.Range("J34").Value = 13944.5
$ws.Range("K34").Value = 2353.0386
$ws.Range("L34").Value = 13944.5
$ws.Range("M34").Value = -2151.0386
$ws.Range("N34").Value = -14348.5
$ws.Range("H132").Value = 76925110
$ws.Range("J132").Value = 1724.5
$ws.Range("L132").Value = 5173.5
$ws.Range("N132").Value = -10233.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 10998.5
$ws.Range("I58").Value = 10998
$ws.Range("K58").Value = 32994
$ws.Range("M58").Value = -32866
$ws.Range("H107").Value = 1345.25
$ws.Range("J107").Value = 1624.3889
$ws.Range("L107").Value = 4873.1667
$ws.Range("N107").Value = -8713.1667
$ws.Range("H122").Value = 1429.1765
$ws.Range("I122").Value = 758.4
$ws.Range("K122").Value = 6825.599999999999
$ws.Range("M122").Value = -4375.599999999999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 56725.156
$ws.Range("I113").Value = 74412.92999999999
$ws.Range("K113").Value = 74412.92999999999
$ws.Range("M113").Value = -72242.92999999999
$ws.Range("H123").Value = 63642.145
$ws.Range("J123").Value = 65082.5
$ws.Range("L123").Value = 65082.5
$ws.Range("N123").Value = -69982.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H92").Value = 45950
$ws.Range("J92").Value = 45950
$ws.Range("L92").Value = 45950
$ws.Range("M92").Value = 0
$ws.Range("N92").Value = -50942
$ws.Range("H132").Value = 120002000
$ws.Range("I132").Value = 240000000
$ws.Range("J132").Value = 3998
$ws.Range("K132").Value = 720000000
$ws.Range("L132").Value = 11994
$ws.Range("M132").Value = -719997470
$ws.Range("N132").Value = -17054
$ws.Range("H136").Value = 2902.9443
$ws.Range("J136").Value = 2998.7856
$ws.Range("L136").Value = 8996.356800000001
$ws.Range("N136").Value = -14096.3568
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4000
$ws.Range("I81").Value = 4000
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 8000
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -6939
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 4000
$ws.Range("I84").Value = 4000
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 40000
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -34696
$ws.Range("N84").Value = ""
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = ""
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = ""
$ws.Range("H107").Value = 1722.5
$ws.Range("I107").Value = 144.5
$ws.Range("K107").Value = 433.5
$ws.Range("M107").Value = 1486.5
$ws.Range("H122").Value = 5469.3
$ws.Range("I122").Value = 5374.1875
$ws.Range("K122").Value = 16122.5625
$ws.Range("M122").Value = -13672.5625
